$wb = $excel.ActiveWorkbook

# Sheet ALC, row 107 (anchor G=27766)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 13).ClearContents()

# Sheet ALC, row 113 (anchor G=27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 2633
$ws.Cells.Item(113, 9).Value = 3000
$ws.Cells.Item(113, 10).Value = 1899
$ws.Cells.Item(113, 11).Value = 3000
$ws.Cells.Item(113, 12).Value = 1899
$ws.Cells.Item(113, 13).Value = 254
$ws.Cells.Item(113, 14).Value = -8407

# Sheet ALC, row 125 (anchor G=36228)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 3864.75
$ws.Cells.Item(125, 9).Value = 987
$ws.Cells.Item(125, 11).Value = 8883
$ws.Cells.Item(125, 13).Value = -6423

# Sheet ALC, row 137 (anchor G=44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3530
$ws.Cells.Item(137, 10).Value = 4717.222
$ws.Cells.Item(137, 12).Value = 14151.666
$ws.Cells.Item(137, 14).Value = -19251.666

# Sheet ALC, row 138 (anchor G=44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4190.5347
$ws.Cells.Item(138, 9).Value = 1611.2858
$ws.Cells.Item(138, 10).Value = 4692.0557
$ws.Cells.Item(138, 11).Value = 4833.857400000001
$ws.Cells.Item(138, 12).Value = 14076.1671
$ws.Cells.Item(138, 13).Value = 306.1425999999992
$ws.Cells.Item(138, 14).Value = -24356.1671

# Sheet ARM, row 45 (anchor G=27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1635.1666
$ws.Cells.Item(45, 9).Value = 1615.25
$ws.Cells.Item(45, 10).Value = 1675
$ws.Cells.Item(45, 11).Value = 1615.25
$ws.Cells.Item(45, 12).Value = 1675
$ws.Cells.Item(45, 13).Value = -1238.25
$ws.Cells.Item(45, 14).Value = -2429

# Sheet ARM, row 61 (anchor G=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4550.5
$ws.Cells.Item(61, 10).Value = 4564.3335
$ws.Cells.Item(61, 12).Value = 4564.3335
$ws.Cells.Item(61, 14).Value = -4988.3335

# Sheet ARM, row 132 (anchor G=43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1451.5454
$ws.Cells.Item(132, 9).Value = 1451.5454
$ws.Cells.Item(132, 11).Value = 4354.6362
$ws.Cells.Item(132, 13).Value = -1824.6362

# Sheet ARM, row 136 (anchor G=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 4550.5
$ws.Cells.Item(136, 10).Value = 4564.3335
$ws.Cells.Item(136, 12).Value = 13693.0005
$ws.Cells.Item(136, 14).Value = -18793.0005

# Sheet BSM, row 86 (anchor G=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1720
$ws.Cells.Item(86, 9).Value = 1720
$ws.Cells.Item(86, 11).Value = 1720
$ws.Cells.Item(86, 13).Value = -597

# Sheet BSM, row 89 (anchor G=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 1720
$ws.Cells.Item(89, 9).Value = 1720
$ws.Cells.Item(89, 11).Value = 8600
$ws.Cells.Item(89, 13).Value = -2984

# Sheet BSM, row 99 (anchor G=19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2830.652
$ws.Cells.Item(99, 9).Value = 2507.1177
$ws.Cells.Item(99, 10).Value = 3747.3333
$ws.Cells.Item(99, 11).Value = 2507.1177
$ws.Cells.Item(99, 12).Value = 3747.3333
$ws.Cells.Item(99, 13).Value = -1009.1177
$ws.Cells.Item(99, 14).Value = -6743.3333

# Sheet BSM, row 107 (anchor G=27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 665.3889
$ws.Cells.Item(107, 9).Value = 665.3889
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 665.3889
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).ClearContents()
$ws.Cells.Item(107, 14).Value = 1254.6111

# Sheet BSM, row 134 (anchor G=43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2916.3333
$ws.Cells.Item(134, 9).Value = 2610.889
$ws.Cells.Item(134, 11).Value = 7832.667
$ws.Cells.Item(134, 13).Value = -5297.667

# Sheet CRP, row 7 (anchor G=5361)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 200.78947
$ws.Cells.Item(7, 10).Value = 499.5
$ws.Cells.Item(7, 12).Value = 499.5
$ws.Cells.Item(7, 14).Value = -725.5

# Sheet CRP, row 22 (anchor G=5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 389.4
$ws.Cells.Item(22, 9).Value = 115.666664
$ws.Cells.Item(22, 10).Value = 800
$ws.Cells.Item(22, 11).Value = 115.666664
$ws.Cells.Item(22, 12).Value = 800
$ws.Cells.Item(22, 13).Value = 234.333336
$ws.Cells.Item(22, 14).Value = -1500

# Sheet CRP, row 58 (anchor G=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3477.12
$ws.Cells.Item(58, 9).Value = 1699.3077
$ws.Cells.Item(58, 10).Value = 5403.0835
$ws.Cells.Item(58, 11).Value = 1699.3077
$ws.Cells.Item(58, 12).Value = 5403.0835
$ws.Cells.Item(58, 13).Value = -1496.3077
$ws.Cells.Item(58, 14).Value = -5809.0835

# Sheet CRP, row 99 (anchor G=36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 11316.361
$ws.Cells.Item(99, 10).Value = 13305
$ws.Cells.Item(99, 12).Value = 13305
$ws.Cells.Item(99, 14).Value = -16301

# Sheet CRP, row 126 (anchor G=36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 11316.361
$ws.Cells.Item(126, 10).Value = 13305
$ws.Cells.Item(126, 12).Value = 39915
$ws.Cells.Item(126, 14).Value = -44855

# Sheet CRP, row 132 (anchor G=44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2316.9443
$ws.Cells.Item(132, 9).Value = 2208
$ws.Cells.Item(132, 10).Value = 2698.25
$ws.Cells.Item(132, 11).Value = 6624
$ws.Cells.Item(132, 12).Value = 8094.75
$ws.Cells.Item(132, 13).Value = -4094
$ws.Cells.Item(132, 14).Value = -13154.75

# Sheet CRP, row 134 (anchor G=44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2570.7144
$ws.Cells.Item(134, 9).Value = 1918.1364
$ws.Cells.Item(134, 11).Value = 5754.4092
$ws.Cells.Item(134, 13).Value = -3219.4092

# Sheet CRP, row 136 (anchor G=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 3477.12
$ws.Cells.Item(136, 9).Value = 1699.3077
$ws.Cells.Item(136, 10).Value = 5403.0835
$ws.Cells.Item(136, 11).Value = 5097.9231
$ws.Cells.Item(136, 12).Value = 16209.2505
$ws.Cells.Item(136, 13).Value = -2547.9231
$ws.Cells.Item(136, 14).Value = -21309.2505

# Sheet CUL, row 5 (anchor G=43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 497.92856
$ws.Cells.Item(5, 9).Value = 517.1667
$ws.Cells.Item(5, 10).Value = 483.5
$ws.Cells.Item(5, 11).Value = 1551.5001
$ws.Cells.Item(5, 12).Value = 1450.5
$ws.Cells.Item(5, 13).Value = -1439.5001
$ws.Cells.Item(5, 14).Value = -1674.5

# Sheet CUL, row 128 (anchor G=41814)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(128, 8).Value = 3979889.2
$ws.Cells.Item(128, 9).Value = 3979889.2
$ws.Cells.Item(128, 11).Value = 11939667.6
$ws.Cells.Item(128, 13).Value = -11934687.6

# Sheet CUL, row 135 (anchor G=43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 497.92856
$ws.Cells.Item(135, 9).Value = 517.1667
$ws.Cells.Item(135, 10).Value = 483.5
$ws.Cells.Item(135, 11).Value = 4654.5003
$ws.Cells.Item(135, 12).Value = 4351.5
$ws.Cells.Item(135, 13).Value = -2119.5003
$ws.Cells.Item(135, 14).Value = -9421.5

# Sheet GSM, row 122 (anchor G=36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 103947.8
$ws.Cells.Item(122, 9).Value = 3176.8
$ws.Cells.Item(122, 10).Value = 204718.8
$ws.Cells.Item(122, 11).Value = 9530.400000000001
$ws.Cells.Item(122, 12).Value = 614156.3999999999
$ws.Cells.Item(122, 13).Value = -7080.400000000001
$ws.Cells.Item(122, 14).Value = -619056.3999999999

# Sheet GSM, row 126 (anchor G=36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 4334.6665
$ws.Cells.Item(126, 9).Value = 3499.5
$ws.Cells.Item(126, 10).Value = 5002.8
$ws.Cells.Item(126, 11).Value = 10498.5
$ws.Cells.Item(126, 12).Value = 15008.4
$ws.Cells.Item(126, 13).Value = -8028.5
$ws.Cells.Item(126, 14).Value = -19948.4

# Sheet GSM, row 132 (anchor G=44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3149.4443
$ws.Cells.Item(132, 9).Value = 1779.5714
$ws.Cells.Item(132, 10).Value = 7944
$ws.Cells.Item(132, 11).Value = 5338.7142
$ws.Cells.Item(132, 12).Value = 23832
$ws.Cells.Item(132, 13).Value = -2808.7142
$ws.Cells.Item(132, 14).Value = -28892

# Sheet LTW, row 16 (anchor G=5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1096.6
$ws.Cells.Item(16, 9).Value = 1159.0714
$ws.Cells.Item(16, 11).Value = 1159.0714
$ws.Cells.Item(16, 13).Value = -989.0714

# Sheet LTW, row 40 (anchor G=36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2236.1428
$ws.Cells.Item(40, 9).Value = 2236.1428
$ws.Cells.Item(40, 11).Value = 2236.1428
$ws.Cells.Item(40, 13).Value = -2100.1428

# Sheet LTW, row 82 (anchor G=12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2930.0688
$ws.Cells.Item(82, 9).Value = 3177.238
$ws.Cells.Item(82, 10).Value = 2281.25
$ws.Cells.Item(82, 11).Value = 3177.238
$ws.Cells.Item(82, 12).Value = 2281.25
$ws.Cells.Item(82, 13).Value = -2816.238
$ws.Cells.Item(82, 14).Value = -3003.25

# Sheet LTW, row 85 (anchor G=12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 2930.0688
$ws.Cells.Item(85, 9).Value = 3177.238
$ws.Cells.Item(85, 10).Value = 2281.25
$ws.Cells.Item(85, 11).Value = 3177.238
$ws.Cells.Item(85, 12).Value = 2281.25
$ws.Cells.Item(85, 13).Value = -1929.238
$ws.Cells.Item(85, 14).Value = -4777.25

# Sheet LTW, row 93 (anchor G=19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1317.3334
$ws.Cells.Item(93, 9).Value = 1313.2
$ws.Cells.Item(93, 10).Value = 1400
$ws.Cells.Item(93, 11).Value = 1313.2
$ws.Cells.Item(93, 12).Value = 1400
$ws.Cells.Item(93, 13).Value = -65.20000000000005
$ws.Cells.Item(93, 14).Value = -3896

# Sheet LTW, row 122 (anchor G=36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 7426.857
$ws.Cells.Item(122, 9).Value = 7331.3335
$ws.Cells.Item(122, 10).Value = 8000
$ws.Cells.Item(122, 11).Value = 21994.0005
$ws.Cells.Item(122, 12).Value = 24000
$ws.Cells.Item(122, 13).Value = -19544.0005
$ws.Cells.Item(122, 14).Value = -28900

# Sheet WVR, row 60 (anchor G=10892)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(60, 8).Value = 38993.5
$ws.Cells.Item(60, 9).Value = 38991
$ws.Cells.Item(60, 11).Value = 38991
$ws.Cells.Item(60, 13).Value = -38169

# Sheet WVR, row 122 (anchor G=36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1103.8
$ws.Cells.Item(122, 9).Value = 1102.25
$ws.Cells.Item(122, 11).Value = 3306.75
$ws.Cells.Item(122, 13).Value = -856.75

# Sheet WVR, row 126 (anchor G=36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2025.3334
$ws.Cells.Item(126, 9).Value = 1366.5555
$ws.Cells.Item(126, 11).Value = 4099.666499999999
$ws.Cells.Item(126, 13).Value = -1629.666499999999

# Sheet WVR, row 132 (anchor G=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1255.4445
$ws.Cells.Item(132, 9).Value = 1129.0416
$ws.Cells.Item(132, 11).Value = 3387.1248
$ws.Cells.Item(132, 13).Value = -857.1248000000001
